# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.672.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "'2.289.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'96.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").Value = "'266.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "'45.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "'0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "'7.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'2.631.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'15.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "'0.844"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "'2.288.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "'43.605.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'71.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "'2.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.09%  "
$ws.Range("D23").Value = "'232.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "'9.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.93%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'2.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "'11.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "'39.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'175.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'21.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").Value = "'0.0883"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "'5.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").Value = "'0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").Value = "'4.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "'12.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +5.81%  "
$ws.Range("D44").Value = "'64.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.03%  "
$ws.Range("D45").Value = "'8.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").Value = "'5.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D48").Value = "'97.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").Value = "'2.510.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  +0.25%  "
